$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- K2 header: "alternatief aliexpress" -> "alternatief " ---
$ws.Range("K2").Value = "alternatief "

# --- Make room for two new rows: one at 23 (koelvin) and one at 25
#     (which becomes the new thermal-pad row once row 24 already holds
#     the re-purposed "thermische pad" data). This pushes the old PCB
#     row (24 -> 26) and the totals row (26 -> 28) down, matching the
#     target layout (gap at 23-row, gap at 25-row kept empty, data rows
#     at 23, 24, 26, 28). ---
$ws.Rows("23:23").Insert()
$ws.Rows("25:25").Insert()

# --- Row 23 (new): koelvin / heatsink alternative from TME ---
$ws.Range("B23").Value = "koelvin"
$ws.Range("C23").Value = " L: 50mm; W: 50mm; H: 50mm; 4.05K/W"
$ws.Range("D23").Value = "ok"
$ws.Range("E23").Value = "TME"
$ws.Range("F23").Value = "https://www.tme.eu/be/en/details/icks50x50x50/heatsinks-for-led/fischer-elektronik/?utm_source=octopart.com&utm_medium=cpc&utm_campaign=compare-2022-12"
$ws.Hyperlinks.Add($ws.Range("F23"), "https://www.tme.eu/be/en/details/icks50x50x50/heatsinks-for-led/fischer-elektronik/?utm_source=octopart.com&utm_medium=cpc&utm_campaign=compare-2022-12")
$ws.Range("F23").Style = "Hyperlink"
$ws.Range("G23").Value = 1
$ws.Range("H23").Value = 12.24
$ws.Range("I23").Formula = "=G23*H23"

# --- Row 24 (re-purposed from "PCB" to "thermische pad") ---
$ws.Range("B24").Value = "thermische pad"
$ws.Range("C24").Value = " SILICONE, 150X0.5MM"
$ws.Range("D24").Value = "ok"
$ws.Range("E24").Value = "farnell"
$ws.Range("F24").Value = "https://be.farnell.com/multicomp-pro/mpgcs-030-150-0-5a/thermal-pad-silicone-150x0-5mm/dp/3267479?st=thermal%20pad%20scilicon"
$ws.Hyperlinks.Add($ws.Range("F24"), "https://be.farnell.com/multicomp-pro/mpgcs-030-150-0-5a/thermal-pad-silicone-150x0-5mm/dp/3267479?st=thermal%20pad%20scilicon")
$ws.Range("F24").Style = "Hyperlink"
$ws.Range("G24").Value = 1
$ws.Range("H24").Value = 13.1
$ws.Range("I24").Formula = "=G24*H24"
$ws.Range("K24").Value = "prijs gedeeld./4"
$ws.Range("L24").Value = 3.275

# --- Row 26 (the original "PCB" row, shifted down here) ---
$ws.Range("B26").Value = "PCB"
$ws.Range("C26").Value = "FR-4 print 1.6mm dikte +stencil"
$ws.Range("E26").Value = "jlcpcb"
$ws.Range("G26").Value = 1
$ws.Range("H26").Value = 13.27
$ws.Range("I26").Formula = "=G26*H26"

# --- Row 28 (totals, was row 26) now also includes the two new rows ---
$ws.Range("B28").Value = "totaal prijs"
$ws.Range("I28").Formula = "=I3+I4+I5+I6+I7+I8+I9+I10+I11+I12+I13+I14+I15+I16+I17+I18+I19+I20+I21+I22+I23+I24+I26"
$ws.Range("L28").Formula = "=I3+I4+I5+I6+L7+I8+I9+I10+I11+I12+I13+I14+I15+I16+I17+I18+I19+I20+I21+I22+I23+L24+I26"

# --- View state: zoom 90%, selection on P27 ---
$ws.Activate()
$excel.ActiveWindow.Zoom = 90
$ws.Range("P27").Select()
